$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'Given  While 1 filling the form, navigate to "Testzen Labs Form" to proceed with registration.'
$ws.Range("C2").Select()
